# MAPPING-COLS-ERD.xlsx - "update README - maps viz"
#
# The sheet has a helper column C that concatenates column name (A) +
# datatype (B) via a shared CONCAT formula. The table-header rows (A35,
# A45, A56, A68, A75 - the bold dataset-name rows) never had a B value,
# so their C cell only ever produced a single trailing-space string that
# wasn't meaningful. This edit clears those five now-unwanted C cells.
#
# It also updates the sheet's view state (zoom level and the active
# selection/scroll position) to reflect where the author was last
# looking in the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stray "<name> " helper-formula cells next to the dataset
# header rows.
$ws.Range("C35").ClearContents()
$ws.Range("C45").ClearContents()
$ws.Range("C56").ClearContents()
$ws.Range("C68").ClearContents()
$ws.Range("C75").ClearContents()

# Update the view: zoom out a bit and scroll/select near the bottom of
# the sheet (last data rows, C76:C77).
$excel.ActiveWindow.Zoom = 73
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 53
$ws.Range("C76:C77").Select()
